$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A handful of the new Price values are plain decimal numerals (e.g. "607.60",
# "0.142", "1.00"). Assigning those strings straight to Range.Value would let
# Excels smart-typing reinterpret them as numbers (and drop the significant
# trailing zeros), so the cells are switched to Text format first. The format
# is switched back to the workbooks Normal style right after the value is
# written, so the text is preserved but no lingering number-format override is
# left behind on the cell.
$textForceCells = @("D5", "D6", "D10", "D11", "D13", "D15", "D19", "D20", "D21", "D22", "D23", "D24", "D27", "D28", "D29", "D31", "D32", "D33", "D34", "D35", "D36", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values cell by cell, in sheet order.
$ws.Range("D2").Value = "66.830.00"
$ws.Range("E2").Value = "  -1.48%  "
$ws.Range("D3").Value = "3.508.76"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "607.60"
$ws.Range("E5").Value = "  +0.15%  "
$ws.Range("D6").Value = "147.39"
$ws.Range("E6").Value = "  -2.31%  "
$ws.Range("D7").Value = "3.505.63"
$ws.Range("E7").Value = "  +0.21%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  -1.93%  "
$ws.Range("D10").Value = "0.142"
$ws.Range("E10").Value = "  -1.27%  "
$ws.Range("D11").Value = "7.99"
$ws.Range("E11").Value = "  +5.57%  "
$ws.Range("E12").Value = "  -2.11%  "
$ws.Range("D13").Value = "0.0000217"
$ws.Range("E13").Value = "  +0.67%  "
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "4.108.73"
$ws.Range("E14").Value = "  +0.32%  "
$ws.Range("B15").Value = "Avalanche"
$ws.Range("C15").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D15").Value = "31.96"
$ws.Range("E15").Value = "  -0.29%  "
$ws.Range("D16").Value = "3.516.92"
$ws.Range("E16").Value = "  +0.13%  "
$ws.Range("D17").Value = "66.988.64"
$ws.Range("E17").Value = "  -1.32%  "
$ws.Range("E18").Value = "  -0.32%  "
$ws.Range("D19").Value = "10.75"
$ws.Range("E19").Value = "  +8.22%  "
$ws.Range("D20").Value = "6.45"
$ws.Range("E20").Value = "  -0.64%  "
$ws.Range("D21").Value = "15.32"
$ws.Range("E21").Value = "  -0.45%  "
$ws.Range("D22").Value = "437.37"
$ws.Range("E22").Value = "  -1.98%  "
$ws.Range("D23").Value = "0.608"
$ws.Range("E23").Value = "  -2.70%  "
$ws.Range("D24").Value = "79.60"
$ws.Range("E24").Value = "  +0.53%  "
$ws.Range("D25").Value = "3.654.16"
$ws.Range("E25").Value = "  +0.36%  "
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("D27").Value = "0.0000120"
$ws.Range("E27").Value = "  -5.16%  "
$ws.Range("D28").Value = "9.76"
$ws.Range("E28").Value = "  -1.92%  "
$ws.Range("D29").Value = "8.22"
$ws.Range("E29").Value = "  -5.20%  "
$ws.Range("E30").Value = "  +0.55%  "
$ws.Range("D31").Value = "1.60"
$ws.Range("E31").Value = "  -2.55%  "
$ws.Range("D32").Value = "0.167"
$ws.Range("E32").Value = "  -2.23%  "
$ws.Range("D33").Value = "0.999"
$ws.Range("E33").Value = "  -0.18%  "
$ws.Range("D34").Value = "25.51"
$ws.Range("E34").Value = "  -0.37%  "
$ws.Range("D35").Value = "5.94"
$ws.Range("E35").Value = "  -3.28%  "
$ws.Range("D36").Value = "1.80"
$ws.Range("E36").Value = "  -2.52%  "
$ws.Range("E37").Value = "  +0.71%  "
$ws.Range("E38").Value = "  +0.02%  "
$ws.Range("D39").Value = "1.00"
$ws.Range("E39").Value = "  +0.03%  "
$ws.Range("D40").Value = "175.52"
$ws.Range("E40").Value = "  -0.88%  "
$ws.Range("D41").Value = "0.0892"
$ws.Range("E41").Value = "  -0.55%  "
$ws.Range("D42").Value = "5.40"
$ws.Range("E42").Value = "  -0.48%  "
$ws.Range("D43").Value = "2.05"
$ws.Range("E43").Value = "  -11.74%  "
$ws.Range("D44").Value = "0.894"
$ws.Range("E44").Value = "  +0.09%  "
$ws.Range("D45").Value = "46.11"
$ws.Range("E45").Value = "  -1.06%  "
$ws.Range("E46").Value = "  -7.55%  "
$ws.Range("D47").Value = "1.25"
$ws.Range("E47").Value = "  -3.36%  "
$ws.Range("D48").Value = "7.45"
$ws.Range("E48").Value = "  -2.02%  "
$ws.Range("D49").Value = "2.43"
$ws.Range("E49").Value = "  -3.80%  "
$ws.Range("D50").Value = "0.991"
$ws.Range("E50").Value = "  +0.12%  "
$ws.Range("D51").Value = "0.247"
$ws.Range("E51").Value = "  -2.00%  "

# Restore the Normal style on the text-forced cells so their formatting matches
# the rest of the sheet (only the cell content changed).
foreach ($addr in $textForceCells) {
    $ws.Range($addr).Style = "Normal"
}
